$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '57.857.93'
$ws.Range("E2").Value = '  -2.35%  '
# Row 3
$ws.Range("D3").Value = '2.451.52'
$ws.Range("E3").Value = '  -3.71%  '
# Row 4
$ws.Range("E4").Value = '  +0.05%  '
# Row 5
$ws.Range("D5").Value = '524.13'
$ws.Range("E5").Value = '  -0.65%  '
# Row 6
$ws.Range("D6").Value = '129.78'
$ws.Range("E6").Value = '  -3.68%  '
# Row 7
$ws.Range("E7").Value = '  -0.06%  '
# Row 8
$ws.Range("E8").Value = '  -0.57%  '
# Row 10
$ws.Range("E10").Value = '  -1.99%  '
# Row 11
$ws.Range("E11").Value = '  -5.03%  '
# Row 12
$ws.Range("E12").Value = '  -4.23%  '
# Row 13
$ws.Range("D13").Value = '2.886.50'
$ws.Range("E13").Value = '  -3.75%  '
# Row 14
$ws.Range("D14").Value = '57.790.75'
$ws.Range("E14").Value = '  -2.34%  '
# Row 15
$ws.Range("E15").Value = '  -4.06%  '
# Row 16
$ws.Range("E16").Value = '  -3.10%  '
# Row 17
$ws.Range("D17").Value = '2.455.07'
$ws.Range("E17").Value = '  -3.52%  '
# Row 18
$ws.Range("D18").Value = '10.39'
$ws.Range("E18").Value = '  -3.35%  '
# Row 19
$ws.Range("E19").Value = '  -2.12%  '
# Row 20
$ws.Range("D20").Value = '312.83'
# Row 21
$ws.Range("E21").Value = '  -1.06%  '
# Row 22
$ws.Range("E22").Value = '  -0.13%  '
# Row 23
$ws.Range("D23").Value = '64.93'
# Row 24
$ws.Range("D24").Value = '0.402'
$ws.Range("E24").Value = '  -2.36%  '
# Row 25
$ws.Range("E25").Value = '  -0.04%  '
# Row 26
$ws.Range("D26").Value = '2.567.72'
$ws.Range("E26").Value = '  -3.55%  '
# Row 27
$ws.Range("E27").Value = '  -3.28%  '
# Row 28
$ws.Range("E28").Value = '  -3.52%  '
# Row 29
$ws.Range("D29").Value = '174.37'
$ws.Range("E29").Value = '  +3.05%  '
# Row 30
$ws.Range("D30").Value = '0.0₃0734'
$ws.Range("E30").Value = '  -3.46%  '
# Row 31
$ws.Range("E31").Value = '  -2.90%  '
# Row 32
$ws.Range("D32").Value = '6.15'
$ws.Range("E32").Value = '  -3.65%  '
# Row 33
$ws.Range("E33").Value = '  -7.80%  '
# Row 34
$ws.Range("E34").Value = '  -0.01%  '
# Row 35
$ws.Range("E35").Value = '  -0.13%  '
# Row 36
$ws.Range("D36").Value = '17.84'
$ws.Range("E36").Value = '  -2.65%  '
# Row 37
$ws.Range("E37").Value = '  -7.59%  '
# Row 38
$ws.Range("E38").Value = '  -5.43%  '
# Row 39
$ws.Range("D39").Value = '36.29'
$ws.Range("E39").Value = '  -1.40%  '
# Row 40
$ws.Range("D40").Value = '0.805'
$ws.Range("E40").Value = '  +2.10%  '
# Row 41
$ws.Range("E41").Value = '  -4.76%  '
# Row 42
$ws.Range("E42").Value = '  -2.95%  '
# Row 43
$ws.Range("E43").Value = '  -3.37%  '
# Row 44
$ws.Range("D44").Value = '4.78'
$ws.Range("E44").Value = '  -6.39%  '
# Row 45
$ws.Range("D45").Value = '257.04'
$ws.Range("E45").Value = '  -8.57%  '
# Row 46
$ws.Range("D46").Value = '123.27'
$ws.Range("E46").Value = '  -8.36%  '
# Row 47
$ws.Range("E47").Value = '  -0.63%  '
# Row 48
$ws.Range("E48").Value = '  -3.22%  '
# Row 49
$ws.Range("D49").Value = "'0.0210"
$ws.Range("E49").Value = '  -3.40%  '
# Row 50
$ws.Range("D50").Value = '16.98'
$ws.Range("E50").Value = '  -5.36%  '
# Row 51
$ws.Range("D51").Value = '16.24'
$ws.Range("E51").Value = '  -5.54%  '
